$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(52, 8).Value = 4787.2
$ws.Cells.Item(52, 10).Value = 4735
$ws.Cells.Item(52, 12).Value = 14205
$ws.Cells.Item(52, 14).Value = -14525
$ws.Cells.Item(112, 8).Value = 2292
$ws.Cells.Item(112, 10).Value = 1996.6666
$ws.Cells.Item(112, 12).Value = 5989.9998
$ws.Cells.Item(112, 14).Value = -8205.9998
$ws.Cells.Item(121, 8).Value = 4367.533
$ws.Cells.Item(121, 10).Value = 4367.533
$ws.Cells.Item(121, 12).Value = 13102.599
$ws.Cells.Item(121, 14).Value = -16596.599
$ws.Cells.Item(127, 8).Value = 787456.0600000001
$ws.Cells.Item(127, 9).Value = 918240.4399999999
$ws.Cells.Item(127, 11).Value = 2754721.32
$ws.Cells.Item(127, 13).Value = -2749761.32
$ws.Cells.Item(129, 8).Value = 1429737.9
$ws.Cells.Item(129, 9).Value = 834485.9399999999
$ws.Cells.Item(129, 11).Value = 2503457.82
$ws.Cells.Item(129, 13).Value = -2498457.82
$ws.Cells.Item(131, 8).Value = 1667864.9
$ws.Cells.Item(131, 9).Value = 1622.5
$ws.Cells.Item(131, 11).Value = 4867.5
$ws.Cells.Item(131, 13).Value = 172.5
$ws.Cells.Item(135, 8).Value = 4761
$ws.Cells.Item(135, 9).Value = 2951.25
$ws.Cells.Item(135, 11).Value = 26561.25
$ws.Cells.Item(135, 13).Value = -24026.25
$ws.Cells.Item(137, 8).Value = 18183182
$ws.Cells.Item(137, 9).Value = 28572914
$ws.Cells.Item(137, 10).Value = 1151.25
$ws.Cells.Item(137, 11).Value = 85718742
$ws.Cells.Item(137, 12).Value = 3453.75
$ws.Cells.Item(137, 13).Value = -85716192
$ws.Cells.Item(137, 14).Value = -8553.75
$ws.Cells.Item(138, 8).Value = 7877.864
$ws.Cells.Item(138, 9).Value = 6839.4
$ws.Cells.Item(138, 11).Value = 20518.2
$ws.Cells.Item(138, 13).Value = -15378.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 875584.0600000001
$ws.Cells.Item(32, 9).Value = 1069612.9
$ws.Cells.Item(32, 11).Value = 1069612.9
$ws.Cells.Item(32, 13).Value = -1069325.9
$ws.Cells.Item(61, 8).Value = 13379947
$ws.Cells.Item(61, 9).Value = 5685337
$ws.Cells.Item(61, 11).Value = 5685337
$ws.Cells.Item(61, 13).Value = -5685125
$ws.Cells.Item(88, 8).Value = 2390.5
$ws.Cells.Item(88, 9).Value = 2601.5
$ws.Cells.Item(88, 10).Value = 2249.8333
$ws.Cells.Item(88, 11).Value = 2601.5
$ws.Cells.Item(88, 12).Value = 2249.8333
$ws.Cells.Item(88, 13).Value = -2195.5
$ws.Cells.Item(88, 14).Value = -3061.8333
$ws.Cells.Item(91, 8).Value = 2390.5
$ws.Cells.Item(91, 9).Value = 2601.5
$ws.Cells.Item(91, 10).Value = 2249.8333
$ws.Cells.Item(91, 11).Value = 2601.5
$ws.Cells.Item(91, 12).Value = 2249.8333
$ws.Cells.Item(91, 13).Value = -1197.5
$ws.Cells.Item(91, 14).Value = -5057.8333
$ws.Cells.Item(97, 8).Value = 523.6842
$ws.Cells.Item(97, 9).Value = 549.4167
$ws.Cells.Item(97, 10).Value = 60.5
$ws.Cells.Item(97, 11).Value = 549.4167
$ws.Cells.Item(97, 12).Value = 60.5
$ws.Cells.Item(97, 13).Value = -53.41669999999999
$ws.Cells.Item(97, 14).Value = -1052.5
$ws.Cells.Item(126, 8).Value = 8509
$ws.Cells.Item(126, 9).Value = 8509
$ws.Cells.Item(126, 11).Value = 25527
$ws.Cells.Item(126, 13).Value = -23057
$ws.Cells.Item(132, 8).Value = 6905.6924
$ws.Cells.Item(132, 9).Value = 4629.6665
$ws.Cells.Item(132, 11).Value = 13888.9995
$ws.Cells.Item(132, 13).Value = -11358.9995
$ws.Cells.Item(136, 8).Value = 13379947
$ws.Cells.Item(136, 9).Value = 5685337
$ws.Cells.Item(136, 11).Value = 17056011
$ws.Cells.Item(136, 13).Value = -17053461

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 65000
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents()  # was -20708
$ws.Cells.Item(86, 8).Value = 4395.5835
$ws.Cells.Item(86, 9).Value = 4343.5
$ws.Cells.Item(86, 11).Value = 4343.5
$ws.Cells.Item(86, 13).Value = -3220.5
$ws.Cells.Item(89, 8).Value = 4395.5835
$ws.Cells.Item(89, 9).Value = 4343.5
$ws.Cells.Item(89, 11).Value = 21717.5
$ws.Cells.Item(89, 13).Value = -16101.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 10590583
$ws.Cells.Item(58, 9).Value = 27781912
$ws.Cells.Item(58, 11).Value = 27781912
$ws.Cells.Item(58, 13).Value = -27781709
$ws.Cells.Item(62, 8).Value = 4323
$ws.Cells.Item(62, 9).Value = 4187.6
$ws.Cells.Item(62, 10).Value = 5000
$ws.Cells.Item(62, 11).Value = 4187.6
$ws.Cells.Item(62, 12).Value = 5000
$ws.Cells.Item(62, 13).Value = -3563.6
$ws.Cells.Item(62, 14).Value = -6248
$ws.Cells.Item(65, 8).Value = 4323
$ws.Cells.Item(65, 9).Value = 4187.6
$ws.Cells.Item(65, 10).Value = 5000
$ws.Cells.Item(65, 11).Value = 20938
$ws.Cells.Item(65, 12).Value = 25000
$ws.Cells.Item(65, 13).Value = -17818
$ws.Cells.Item(65, 14).Value = -31240
$ws.Cells.Item(99, 8).Value = 3622.7222
$ws.Cells.Item(99, 9).Value = 2509.25
$ws.Cells.Item(99, 11).Value = 2509.25
$ws.Cells.Item(99, 13).Value = -1011.25
$ws.Cells.Item(126, 8).Value = 3622.7222
$ws.Cells.Item(126, 9).Value = 2509.25
$ws.Cells.Item(126, 11).Value = 7527.75
$ws.Cells.Item(126, 13).Value = -5057.75
$ws.Cells.Item(134, 8).Value = 4490.2354
$ws.Cells.Item(134, 9).Value = 2202.2
$ws.Cells.Item(134, 11).Value = 6606.599999999999
$ws.Cells.Item(134, 13).Value = -4071.599999999999
$ws.Cells.Item(136, 8).Value = 10590583
$ws.Cells.Item(136, 9).Value = 27781912
$ws.Cells.Item(136, 11).Value = 83345736
$ws.Cells.Item(136, 13).Value = -83343186
$ws.Cells.Item(140, 8).Value = 68995
$ws.Cells.Item(140, 10).Value = 68995
$ws.Cells.Item(140, 12).Value = 68995
$ws.Cells.Item(140, 14).Value = -79355
$ws.Cells.Item(141, 8).Value = 78326
$ws.Cells.Item(141, 10).Value = 78326
$ws.Cells.Item(141, 12).Value = 78326
$ws.Cells.Item(141, 13).Value = -88686

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 4374.4375
$ws.Cells.Item(140, 9).Value = 2789.2
$ws.Cells.Item(140, 11).Value = 8367.599999999999
$ws.Cells.Item(140, 13).Value = -3187.599999999999
$ws.Cells.Item(141, 8).Value = 12483.333
$ws.Cells.Item(141, 9).Value = 1450
$ws.Cells.Item(141, 11).Value = 4350
$ws.Cells.Item(141, 13).Value = 830

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 8428.286
$ws.Cells.Item(80, 9).Value = 7666.3335
$ws.Cells.Item(80, 10).Value = 8999.75
$ws.Cells.Item(80, 11).Value = 7666.3335
$ws.Cells.Item(80, 12).Value = 8999.75
$ws.Cells.Item(80, 13).Value = -6668.3335
$ws.Cells.Item(80, 14).Value = -10995.75
$ws.Cells.Item(83, 8).Value = 8428.286
$ws.Cells.Item(83, 9).Value = 7666.3335
$ws.Cells.Item(83, 10).Value = 8999.75
$ws.Cells.Item(83, 11).Value = 38331.6675
$ws.Cells.Item(83, 12).Value = 44998.75
$ws.Cells.Item(83, 13).Value = -33339.6675
$ws.Cells.Item(83, 14).Value = -54982.75
$ws.Cells.Item(93, 8).Value = 39000
$ws.Cells.Item(93, 10).Value = 39000
$ws.Cells.Item(93, 12).Value = 39000
$ws.Cells.Item(93, 14).Value = -42744
$ws.Cells.Item(122, 8).Value = 58649.668
$ws.Cells.Item(122, 9).Value = 103104.8
$ws.Cells.Item(122, 10).Value = 3080.75
$ws.Cells.Item(122, 11).Value = 309314.4
$ws.Cells.Item(122, 12).Value = 9242.25
$ws.Cells.Item(122, 13).Value = -306864.4
$ws.Cells.Item(122, 14).Value = -14142.25
$ws.Cells.Item(136, 8).Value = 29860.5
$ws.Cells.Item(136, 10).Value = 29269.143
$ws.Cells.Item(136, 12).Value = 87807.429
$ws.Cells.Item(136, 14).Value = -92907.429

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 2968.8462
$ws.Cells.Item(93, 9).Value = 1574.6666
$ws.Cells.Item(93, 11).Value = 1574.6666
$ws.Cells.Item(93, 13).Value = -326.6666
$ws.Cells.Item(122, 8).Value = 4880.125
$ws.Cells.Item(122, 9).Value = 4602.75
$ws.Cells.Item(122, 10).Value = 5712.25
$ws.Cells.Item(122, 11).Value = 13808.25
$ws.Cells.Item(122, 12).Value = 17136.75
$ws.Cells.Item(122, 13).Value = -11358.25
$ws.Cells.Item(122, 14).Value = -22036.75
$ws.Cells.Item(132, 8).Value = 4170266
$ws.Cells.Item(132, 9).Value = 11907435
$ws.Cells.Item(132, 10).Value = 4098.077
$ws.Cells.Item(132, 11).Value = 35722305
$ws.Cells.Item(132, 12).Value = 12294.231
$ws.Cells.Item(132, 13).Value = -35719775
$ws.Cells.Item(132, 14).Value = -17354.231

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4024.625
$ws.Cells.Item(81, 9).Value = 3956.7144
$ws.Cells.Item(81, 10).Value = 4500
$ws.Cells.Item(81, 11).Value = 7913.4288
$ws.Cells.Item(81, 12).Value = 9000
$ws.Cells.Item(81, 13).Value = -6852.4288
$ws.Cells.Item(81, 14).Value = -11122
$ws.Cells.Item(84, 8).Value = 4024.625
$ws.Cells.Item(84, 9).Value = 3956.7144
$ws.Cells.Item(84, 10).Value = 4500
$ws.Cells.Item(84, 11).Value = 39567.144
$ws.Cells.Item(84, 12).Value = 45000
$ws.Cells.Item(84, 13).Value = -34263.144
$ws.Cells.Item(84, 14).Value = -55608
$ws.Cells.Item(126, 8).Value = 1361.5927
$ws.Cells.Item(126, 9).Value = 1126.2858
$ws.Cells.Item(126, 11).Value = 3378.8574
$ws.Cells.Item(126, 13).Value = -908.8574000000003
